# Auto-generated Excel COM-interop script applying the cryptos.xlsx price/volume refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some refreshed Price values (column D) are plain decimals (e.g. "1.00", "609.75").
# Excel's normal text->value coercion would turn those into real numbers and drop the
# literal formatting ("1.00" -> 1). Mark just those cells as Text first so the
# assignments below stick as the exact strings scraped from the source site.
# (Set one cell at a time - this host's Range() does not support comma-separated
# unions/multi-arg corners for property assignment.)
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"

$ws.Range("D2").Value = '64.413.05'
$ws.Range("E2").Value = '  +1.75%  '
$ws.Range("D3").Value = '2.654.40'
$ws.Range("E3").Value = '  -0.59%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '609.75'
$ws.Range("E5").Value = '  -0.61%  '
$ws.Range("D6").Value = '148.48'
$ws.Range("E6").Value = '  +3.54%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '0.592'
$ws.Range("E8").Value = '  +0.77%  '
$ws.Range("E9").Value = '  +2.79%  '
$ws.Range("D10").Value = '0.388'
$ws.Range("E10").Value = '  +7.17%  '
$ws.Range("D11").Value = '5.61'
$ws.Range("E11").Value = '  +0.13%  '
$ws.Range("E12").Value = '  -0.94%  '
$ws.Range("D13").Value = '27.72'
$ws.Range("E13").Value = '  +1.36%  '
$ws.Range("D14").Value = '3.125.19'
$ws.Range("E14").Value = '  -0.83%  '
$ws.Range("D15").Value = '64.226.92'
$ws.Range("E15").Value = '  +1.68%  '
$ws.Range("E16").Value = '  +2.23%  '
$ws.Range("D17").Value = '2.651.88'
$ws.Range("E17").Value = '  -0.78%  '
$ws.Range("D18").Value = '11.88'
$ws.Range("E18").Value = '  +4.05%  '
$ws.Range("D19").Value = '4.60'
$ws.Range("E19").Value = '  +4.28%  '
$ws.Range("D20").Value = '347.17'
$ws.Range("E20").Value = '  +1.50%  '
$ws.Range("D21").Value = '6.92'
$ws.Range("E21").Value = '  +0.68%  '
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("D23").Value = '5.56'
$ws.Range("E23").Value = '  -1.42%  '
$ws.Range("D24").Value = '66.75'
$ws.Range("E24").Value = '  -0.71%  '
$ws.Range("E25").Value = '  +8.16%  '
$ws.Range("D26").Value = '9.44'
$ws.Range("E26").Value = '  +8.54%  '
$ws.Range("D27").Value = '1.71'
$ws.Range("E27").Value = '  +3.39%  '
$ws.Range("D28").Value = '558.53'
$ws.Range("E28").Value = '  +3.29%  '
$ws.Range("D29").Value = '8.17'
$ws.Range("E29").Value = '  +3.15%  '
$ws.Range("D30").Value = '0.998'
$ws.Range("E30").Value = '  -0.23%  '
$ws.Range("E31").Value = '  -2.28%  '
$ws.Range("D32").Value = '2.07'
$ws.Range("E32").Value = '  +0.82%  '
$ws.Range("E33").Value = '  +6.52%  '
$ws.Range("E34").Value = '  -1.20%  '
$ws.Range("E35").Value = '  +3.10%  '
$ws.Range("D36").Value = '169.45'
$ws.Range("E36").Value = '  -1.57%  '
$ws.Range("D37").Value = '0.406'
$ws.Range("E37").Value = '  +0.21%  '
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  +0.17%  '
$ws.Range("D39").Value = '1.95'
$ws.Range("E39").Value = '  +4.97%  '
$ws.Range("D40").Value = '19.38'
$ws.Range("E40").Value = '  +1.00%  '
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("D42").Value = '165.52'
$ws.Range("E42").Value = '  -6.36%  '
$ws.Range("D43").Value = '40.34'
$ws.Range("E43").Value = '  +0.60%  '
$ws.Range("D44").Value = '3.84'
$ws.Range("E44").Value = '  +2.79%  '
$ws.Range("D45").Value = '22.16'
$ws.Range("E45").Value = '  -0.80%  '
$ws.Range("E46").Value = '  -0.16%  '
$ws.Range("D47").Value = '0.630'
$ws.Range("E47").Value = '  -0.66%  '
$ws.Range("D50").Value = '0.0962'
$ws.Range("E50").Value = '  -0.21%  '
$ws.Range("D51").Value = '19.01'
$ws.Range("E51").Value = '  +0.95%  '

# Rows 48 and 49 swapped rank order (dogwifhat <-> VeChain) with refreshed figures.
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").Value = '0.0247'
$ws.Range("E48").Value = '  +3.18%  '
$ws.Range("B49").Value = 'dogwifhat'
$ws.Range("C49").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D49").Value = '2.00'
$ws.Range("E49").Value = '  +14.59%  '
